$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows for the "ECs" sending-cluster block above the existing data
$ws.Rows("2:7").Insert()

# Row 2: ECs -> ECs
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Has2"
$ws.Cells.Item(2,3).Value = "Hmmr"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 0.1403053333333333
$ws.Cells.Item(2,8).Value = 0.420916
$ws.Cells.Item(2,9).Value = 0.002676192217864005
$ws.Cells.Item(2,10).Value = 0.002676192217864005
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 3.457962
$ws.Cells.Item(2,14).Value = 10.373886
$ws.Cells.Item(2,15).Value = 0.07614454177940357
$ws.Cells.Item(2,16).Value = 0.07614454177940358
$ws.Cells.Item(2,17).Value = 0.4851705110640001
$ws.Cells.Item(2,18).Value = 4.366534599576
$ws.Cells.Item(2,19).Value = 0.0002037774301428604
$ws.Cells.Item(2,20).Value = 0.0002037774301428604

# Row 3: ECs -> FAPs
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Has2"
$ws.Cells.Item(3,3).Value = "Hmmr"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 0.1403053333333333
$ws.Cells.Item(3,8).Value = 0.420916
$ws.Cells.Item(3,9).Value = 0.002676192217864005
$ws.Cells.Item(3,10).Value = 0.002676192217864005
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 25.73962066666667
$ws.Cells.Item(3,14).Value = 77.218862
$ws.Cells.Item(3,15).Value = 0.566788073795779
$ws.Cells.Item(3,16).Value = 0.566788073795779
$ws.Cells.Item(3,17).Value = 3.611406057510222
$ws.Cells.Item(3,18).Value = 32.502654517592
$ws.Cells.Item(3,19).Value = 0.001516833832270393
$ws.Cells.Item(3,20).Value = 0.001516833832270393

# Row 4: ECs -> Inflammatory-Mac
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Has2"
$ws.Cells.Item(4,3).Value = "Hmmr"
$ws.Cells.Item(4,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 0.1403053333333333
$ws.Cells.Item(4,8).Value = 0.420916
$ws.Cells.Item(4,9).Value = 0.002676192217864005
$ws.Cells.Item(4,10).Value = 0.002676192217864005
$ws.Cells.Item(4,11).Value = 2
$ws.Cells.Item(4,12).Value = 0.6666666666666666
$ws.Cells.Item(4,13).Value = 0.3208433333333333
$ws.Cells.Item(4,14).Value = 0.96253
$ws.Cells.Item(4,15).Value = 0.007064990476946567
$ws.Cells.Item(4,16).Value = 0.007064990476946567
$ws.Cells.Item(4,17).Value = 0.04501603083111111
$ws.Cells.Item(4,18).Value = 0.40514427748
$ws.Cells.Item(4,19).Value = 0.000018907272533687698856793777
$ws.Cells.Item(4,20).Value = 0.000018907272533687709021189144

# Row 5: ECs -> MuSCs
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Has2"
$ws.Cells.Item(5,3).Value = "Hmmr"
$ws.Cells.Item(5,4).Value = "MuSCs"
$ws.Cells.Item(5,5).Value = 1
$ws.Cells.Item(5,6).Value = 0.3333333333333333
$ws.Cells.Item(5,7).Value = 0.1403053333333333
$ws.Cells.Item(5,8).Value = 0.420916
$ws.Cells.Item(5,9).Value = 0.002676192217864005
$ws.Cells.Item(5,10).Value = 0.002676192217864005
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 12.82934566666667
$ws.Cells.Item(5,14).Value = 38.48803700000001
$ws.Cells.Item(5,15).Value = 0.282503002380567
$ws.Cells.Item(5,16).Value = 0.282503002380567
$ws.Cells.Item(5,17).Value = 1.800025620210222
$ws.Cells.Item(5,18).Value = 16.200230581892
$ws.Cells.Item(5,19).Value = 0.0007560323364940899
$ws.Cells.Item(5,20).Value = 0.0007560323364940899

# Row 6: ECs -> Neutrophils
$ws.Cells.Item(6,1).Value = "ECs"
$ws.Cells.Item(6,2).Value = "Has2"
$ws.Cells.Item(6,3).Value = "Hmmr"
$ws.Cells.Item(6,4).Value = "Neutrophils"
$ws.Cells.Item(6,5).Value = 1
$ws.Cells.Item(6,6).Value = 0.3333333333333333
$ws.Cells.Item(6,7).Value = 0.1403053333333333
$ws.Cells.Item(6,8).Value = 0.420916
$ws.Cells.Item(6,9).Value = 0.002676192217864005
$ws.Cells.Item(6,10).Value = 0.002676192217864005
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 1.320849666666667
$ws.Cells.Item(6,14).Value = 3.962549
$ws.Cells.Item(6,15).Value = 0.02908519313624941
$ws.Cells.Item(6,16).Value = 0.02908519313624941
$ws.Cells.Item(6,17).Value = 0.1853222527648889
$ws.Cells.Item(6,18).Value = 1.667900274884
$ws.Cells.Item(6,19).Value = 0.000077837567526302218904113439
$ws.Cells.Item(6,20).Value = 0.000077837567526302232456640595

# Row 7: ECs -> Resolving-Mac
$ws.Cells.Item(7,1).Value = "ECs"
$ws.Cells.Item(7,2).Value = "Has2"
$ws.Cells.Item(7,3).Value = "Hmmr"
$ws.Cells.Item(7,4).Value = "Resolving-Mac"
$ws.Cells.Item(7,5).Value = 1
$ws.Cells.Item(7,6).Value = 0.3333333333333333
$ws.Cells.Item(7,7).Value = 0.1403053333333333
$ws.Cells.Item(7,8).Value = 0.420916
$ws.Cells.Item(7,9).Value = 0.002676192217864005
$ws.Cells.Item(7,10).Value = 0.002676192217864005
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 1.744509
$ws.Cells.Item(7,14).Value = 5.233527
$ws.Cells.Item(7,15).Value = 0.03841419843105434
$ws.Cells.Item(7,16).Value = 0.03841419843105434
$ws.Cells.Item(7,17).Value = 0.244763916748
$ws.Cells.Item(7,18).Value = 2.202875250732
$ws.Cells.Item(7,19).Value = 0.0001028037788966713
$ws.Cells.Item(7,20).Value = 0.0001028037788966713

# Row 8: FAPs -> ECs
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Has2"
$ws.Cells.Item(8,3).Value = "Hmmr"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 48.25514733333333
$ws.Cells.Item(8,8).Value = 144.765442
$ws.Cells.Item(8,9).Value = 0.9204215313650299
$ws.Cells.Item(8,10).Value = 0.9204215313650299
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 3.457962
$ws.Cells.Item(8,14).Value = 10.373886
$ws.Cells.Item(8,15).Value = 0.07614454177940357
$ws.Cells.Item(8,16).Value = 0.07614454177940358
$ws.Cells.Item(8,17).Value = 166.864465783068
$ws.Cells.Item(8,18).Value = 1501.780192047612
$ws.Cells.Item(8,19).Value = 0.07008507574968713
$ws.Cells.Item(8,20).Value = 0.07008507574968714

# Row 9: FAPs -> FAPs
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Has2"
$ws.Cells.Item(9,3).Value = "Hmmr"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 48.25514733333333
$ws.Cells.Item(9,8).Value = 144.765442
$ws.Cells.Item(9,9).Value = 0.9204215313650299
$ws.Cells.Item(9,10).Value = 0.9204215313650299
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 25.73962066666667
$ws.Cells.Item(9,14).Value = 77.218862
$ws.Cells.Item(9,15).Value = 0.566788073795779
$ws.Cells.Item(9,16).Value = 0.566788073795779
$ws.Cells.Item(9,17).Value = 1242.069187574112
$ws.Cells.Item(9,18).Value = 11178.622688167
$ws.Cells.Item(9,19).Value = 0.5216839468425465
$ws.Cells.Item(9,20).Value = 0.5216839468425465

# Row 10: FAPs -> Inflammatory-Mac
$ws.Cells.Item(10,1).Value = "FAPs"
$ws.Cells.Item(10,2).Value = "Has2"
$ws.Cells.Item(10,3).Value = "Hmmr"
$ws.Cells.Item(10,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 48.25514733333333
$ws.Cells.Item(10,8).Value = 144.765442
$ws.Cells.Item(10,9).Value = 0.9204215313650299
$ws.Cells.Item(10,10).Value = 0.9204215313650299
$ws.Cells.Item(10,11).Value = 2
$ws.Cells.Item(10,12).Value = 0.6666666666666666
$ws.Cells.Item(10,13).Value = 0.3208433333333333
$ws.Cells.Item(10,14).Value = 0.96253
$ws.Cells.Item(10,15).Value = 0.007064990476946567
$ws.Cells.Item(10,16).Value = 0.007064990476946567
$ws.Cells.Item(10,17).Value = 15.48234232091778
$ws.Cells.Item(10,18).Value = 139.34108088826
$ws.Cells.Item(10,19).Value = 0.006502769353870512
$ws.Cells.Item(10,20).Value = 0.006502769353870513

# Row 11: FAPs -> MuSCs
$ws.Cells.Item(11,1).Value = "FAPs"
$ws.Cells.Item(11,2).Value = "Has2"
$ws.Cells.Item(11,3).Value = "Hmmr"
$ws.Cells.Item(11,4).Value = "MuSCs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 48.25514733333333
$ws.Cells.Item(11,8).Value = 144.765442
$ws.Cells.Item(11,9).Value = 0.9204215313650299
$ws.Cells.Item(11,10).Value = 0.9204215313650299
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 12.82934566666667
$ws.Cells.Item(11,14).Value = 38.48803700000001
$ws.Cells.Item(11,15).Value = 0.282503002380567
$ws.Cells.Item(11,16).Value = 0.282503002380567
$ws.Cells.Item(11,17).Value = 619.0819653352617
$ws.Cells.Item(11,18).Value = 5571.737688017355
$ws.Cells.Item(11,19).Value = 0.2600218460663402
$ws.Cells.Item(11,20).Value = 0.2600218460663402

# Row 12: FAPs -> Neutrophils
$ws.Cells.Item(12,1).Value = "FAPs"
$ws.Cells.Item(12,2).Value = "Has2"
$ws.Cells.Item(12,3).Value = "Hmmr"
$ws.Cells.Item(12,4).Value = "Neutrophils"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 48.25514733333333
$ws.Cells.Item(12,8).Value = 144.765442
$ws.Cells.Item(12,9).Value = 0.9204215313650299
$ws.Cells.Item(12,10).Value = 0.9204215313650299
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 1.320849666666667
$ws.Cells.Item(12,14).Value = 3.962549
$ws.Cells.Item(12,15).Value = 0.02908519313624941
$ws.Cells.Item(12,16).Value = 0.02908519313624941
$ws.Cells.Item(12,17).Value = 63.73779527018422
$ws.Cells.Item(12,18).Value = 573.640157431658
$ws.Cells.Item(12,19).Value = 0.02677063800651434
$ws.Cells.Item(12,20).Value = 0.02677063800651434

# Row 13: FAPs -> Resolving-Mac
$ws.Cells.Item(13,1).Value = "FAPs"
$ws.Cells.Item(13,2).Value = "Has2"
$ws.Cells.Item(13,3).Value = "Hmmr"
$ws.Cells.Item(13,4).Value = "Resolving-Mac"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 48.25514733333333
$ws.Cells.Item(13,8).Value = 144.765442
$ws.Cells.Item(13,9).Value = 0.9204215313650299
$ws.Cells.Item(13,10).Value = 0.9204215313650299
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 1.744509
$ws.Cells.Item(13,14).Value = 5.233527
$ws.Cells.Item(13,15).Value = 0.03841419843105434
$ws.Cells.Item(13,16).Value = 0.03841419843105434
$ws.Cells.Item(13,17).Value = 84.181538819326
$ws.Cells.Item(13,18).Value = 757.633849373934
$ws.Cells.Item(13,19).Value = 0.03535725534607117
$ws.Cells.Item(13,20).Value = 0.03535725534607117

# Row 14: MuSCs -> ECs
$ws.Cells.Item(14,1).Value = "MuSCs"
$ws.Cells.Item(14,2).Value = "Has2"
$ws.Cells.Item(14,3).Value = "Hmmr"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 4.031773
$ws.Cells.Item(14,8).Value = 12.095319
$ws.Cells.Item(14,9).Value = 0.07690227641710612
$ws.Cells.Item(14,10).Value = 0.0769022764171061
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(14,13).Value = 3.457962
$ws.Cells.Item(14,14).Value = 10.373886
$ws.Cells.Item(14,15).Value = 0.07614454177940357
$ws.Cells.Item(14,16).Value = 0.07614454177940358
$ws.Cells.Item(14,17).Value = 13.941717826626
$ws.Cells.Item(14,18).Value = 125.475460439634
$ws.Cells.Item(14,19).Value = 0.005855688599573579
$ws.Cells.Item(14,20).Value = 0.005855688599573579

# Row 15: MuSCs -> FAPs
$ws.Cells.Item(15,1).Value = "MuSCs"
$ws.Cells.Item(15,2).Value = "Has2"
$ws.Cells.Item(15,3).Value = "Hmmr"
$ws.Cells.Item(15,4).Value = "FAPs"
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 4.031773
$ws.Cells.Item(15,8).Value = 12.095319
$ws.Cells.Item(15,9).Value = 0.07690227641710612
$ws.Cells.Item(15,10).Value = 0.0769022764171061
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 25.73962066666667
$ws.Cells.Item(15,14).Value = 77.218862
$ws.Cells.Item(15,15).Value = 0.566788073795779
$ws.Cells.Item(15,16).Value = 0.566788073795779
$ws.Cells.Item(15,17).Value = 103.7763076341087
$ws.Cells.Item(15,18).Value = 933.9867687069781
$ws.Cells.Item(15,19).Value = 0.04358729312096214
$ws.Cells.Item(15,20).Value = 0.04358729312096213

# Row 16: MuSCs -> Inflammatory-Mac
$ws.Cells.Item(16,1).Value = "MuSCs"
$ws.Cells.Item(16,2).Value = "Has2"
$ws.Cells.Item(16,3).Value = "Hmmr"
$ws.Cells.Item(16,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 4.031773
$ws.Cells.Item(16,8).Value = 12.095319
$ws.Cells.Item(16,9).Value = 0.07690227641710612
$ws.Cells.Item(16,10).Value = 0.0769022764171061
$ws.Cells.Item(16,11).Value = 2
$ws.Cells.Item(16,12).Value = 0.6666666666666666
$ws.Cells.Item(16,13).Value = 0.3208433333333333
$ws.Cells.Item(16,14).Value = 0.96253
$ws.Cells.Item(16,15).Value = 0.007064990476946567
$ws.Cells.Item(16,16).Value = 0.007064990476946567
$ws.Cells.Item(16,17).Value = 1.293567488563333
$ws.Cells.Item(16,18).Value = 11.64210739707
$ws.Cells.Item(16,19).Value = 0.0005433138505423672
$ws.Cells.Item(16,20).Value = 0.0005433138505423672

# Row 17: MuSCs -> MuSCs
$ws.Cells.Item(17,1).Value = "MuSCs"
$ws.Cells.Item(17,2).Value = "Has2"
$ws.Cells.Item(17,3).Value = "Hmmr"
$ws.Cells.Item(17,4).Value = "MuSCs"
$ws.Cells.Item(17,5).Value = 3
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(17,7).Value = 4.031773
$ws.Cells.Item(17,8).Value = 12.095319
$ws.Cells.Item(17,9).Value = 0.07690227641710612
$ws.Cells.Item(17,10).Value = 0.0769022764171061
$ws.Cells.Item(17,11).Value = 3
$ws.Cells.Item(17,12).Value = 1
$ws.Cells.Item(17,13).Value = 12.82934566666667
$ws.Cells.Item(17,14).Value = 38.48803700000001
$ws.Cells.Item(17,15).Value = 0.282503002380567
$ws.Cells.Item(17,16).Value = 0.282503002380567
$ws.Cells.Item(17,17).Value = 51.72500946653368
$ws.Cells.Item(17,18).Value = 465.5250851988031
$ws.Cells.Item(17,19).Value = 0.02172512397773275
$ws.Cells.Item(17,20).Value = 0.02172512397773275

# Row 18: MuSCs -> Neutrophils
$ws.Cells.Item(18,1).Value = "MuSCs"
$ws.Cells.Item(18,2).Value = "Has2"
$ws.Cells.Item(18,3).Value = "Hmmr"
$ws.Cells.Item(18,4).Value = "Neutrophils"
$ws.Cells.Item(18,5).Value = 3
$ws.Cells.Item(18,6).Value = 1
$ws.Cells.Item(18,7).Value = 4.031773
$ws.Cells.Item(18,8).Value = 12.095319
$ws.Cells.Item(18,9).Value = 0.07690227641710612
$ws.Cells.Item(18,10).Value = 0.0769022764171061
$ws.Cells.Item(18,11).Value = 3
$ws.Cells.Item(18,12).Value = 1
$ws.Cells.Item(18,13).Value = 1.320849666666667
$ws.Cells.Item(18,14).Value = 3.962549
$ws.Cells.Item(18,15).Value = 0.02908519313624941
$ws.Cells.Item(18,16).Value = 0.02908519313624941
$ws.Cells.Item(18,17).Value = 5.325366023125667
$ws.Cells.Item(18,18).Value = 47.928294208131
$ws.Cells.Item(18,19).Value = 0.00223671756220877
$ws.Cells.Item(18,20).Value = 0.002236717562208769

# Row 19: MuSCs -> Resolving-Mac
$ws.Cells.Item(19,1).Value = "MuSCs"
$ws.Cells.Item(19,2).Value = "Has2"
$ws.Cells.Item(19,3).Value = "Hmmr"
$ws.Cells.Item(19,4).Value = "Resolving-Mac"
$ws.Cells.Item(19,5).Value = 3
$ws.Cells.Item(19,6).Value = 1
$ws.Cells.Item(19,7).Value = 4.031773
$ws.Cells.Item(19,8).Value = 12.095319
$ws.Cells.Item(19,9).Value = 0.07690227641710612
$ws.Cells.Item(19,10).Value = 0.0769022764171061
$ws.Cells.Item(19,11).Value = 3
$ws.Cells.Item(19,12).Value = 1
$ws.Cells.Item(19,13).Value = 1.744509
$ws.Cells.Item(19,14).Value = 5.233527
$ws.Cells.Item(19,15).Value = 0.03841419843105434
$ws.Cells.Item(19,16).Value = 0.03841419843105434
$ws.Cells.Item(19,17).Value = 7.033464284457
$ws.Cells.Item(19,18).Value = 63.30117856011299
$ws.Cells.Item(19,19).Value = 0.002954139306086505
$ws.Cells.Item(19,20).Value = 0.002954139306086505
